$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 384.36365
$ws.Range("I2").Value = 142.5
$ws.Range("K2").Value = 142.5
$ws.Range("M2").Value = -29.5
$ws.Range("H9").Value = 5676.737
$ws.Range("I9").Value = 6690.4375
$ws.Range("J9").Value = 270.33334
$ws.Range("K9").Value = 6690.4375
$ws.Range("L9").Value = 270.33334
$ws.Range("M9").Value = -6521.4375
$ws.Range("N9").Value = -608.33334
$ws.Range("H18").Value = 1110
$ws.Range("I18").Value = 1110
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1110
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -826
$ws.Range("N18").ClearContents()
$ws.Range("H29").Value = 3635.5715
$ws.Range("J29").Value = 6059.75
$ws.Range("L29").Value = 18179.25
$ws.Range("N29").Value = -18741.25
$ws.Range("H32").Value = 7508.857
$ws.Range("I32").Value = 6635
$ws.Range("J32").Value = 7858.4
$ws.Range("K32").Value = 6635
$ws.Range("L32").Value = 7858.4
$ws.Range("M32").Value = -6309
$ws.Range("N32").Value = -8510.4
$ws.Range("H33").Value = 238.77777
$ws.Range("J33").Value = 477.25
$ws.Range("L33").Value = 477.25
$ws.Range("N33").Value = -935.25
$ws.Range("H64").Value = 28575500
$ws.Range("I64").Value = 50003372
$ws.Range("K64").Value = 50003372
$ws.Range("M64").Value = -50003124
$ws.Range("H67").Value = 28575500
$ws.Range("I67").Value = 50003372
$ws.Range("K67").Value = 50003372
$ws.Range("M67").Value = -50002514
$ws.Range("H74").Value = 41673950
$ws.Range("I74").Value = 41673950
$ws.Range("K74").Value = 41673950
$ws.Range("M74").Value = -41673014
$ws.Range("H77").Value = 41673950
$ws.Range("I77").Value = 41673950
$ws.Range("K77").Value = 208369750
$ws.Range("M77").Value = -208365070
$ws.Range("H86").Value = 10159.444
$ws.Range("J86").Value = 3731
$ws.Range("L86").Value = 3731
$ws.Range("N86").Value = -5977
$ws.Range("H88").Value = 52729970
$ws.Range("I88").Value = 111112830
$ws.Range("J88").Value = 8942821
$ws.Range("K88").Value = 111112830
$ws.Range("L88").Value = 8942821
$ws.Range("M88").Value = -111112424
$ws.Range("N88").Value = -8943633
$ws.Range("H89").Value = 10159.444
$ws.Range("J89").Value = 3731
$ws.Range("L89").Value = 18655
$ws.Range("N89").Value = -29887
$ws.Range("H91").Value = 52729970
$ws.Range("I91").Value = 111112830
$ws.Range("J91").Value = 8942821
$ws.Range("K91").Value = 111112830
$ws.Range("L91").Value = 8942821
$ws.Range("M91").Value = -111111426
$ws.Range("N91").Value = -8945629
$ws.Range("H111").Value = 1689.5
$ws.Range("J111").Value = 1350
$ws.Range("L111").Value = 4050
$ws.Range("N111").Value = -10184
$ws.Range("H113").Value = 3995
$ws.Range("I113").Value = 3995
$ws.Range("K113").Value = 3995
$ws.Range("M113").Value = -741
$ws.Range("H116").Value = 3921.3157
$ws.Range("I116").Value = 4015.7693
$ws.Range("K116").Value = 4015.7693
$ws.Range("M116").Value = -573.7692999999999
$ws.Range("H137").Value = 5194.722
$ws.Range("I137").Value = 5172.1665
$ws.Range("K137").Value = 15516.4995
$ws.Range("M137").Value = -12966.4995
$ws.Range("H138").Value = 6973.323
$ws.Range("J138").Value = 8156.022
$ws.Range("L138").Value = 24468.066
$ws.Range("N138").Value = -34748.066
$ws.Range("H141").Value = 6098.1665
$ws.Range("I141").Value = 4172.25
$ws.Range("K141").Value = 12516.75
$ws.Range("M141").Value = -7336.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 921140.75
$ws.Range("J2").Value = 4124.75
$ws.Range("K2").Value = 921140.75
$ws.Range("L2").Value = 4124.75
$ws.Range("M2").Value = -921027.75
$ws.Range("N2").Value = -4350.75
$ws.Range("H5").Value = 748.4
$ws.Range("I5").Value = 645.2941
$ws.Range("J5").Value = 1332.6666
$ws.Range("K5").Value = 645.2941
$ws.Range("L5").Value = 1332.6666
$ws.Range("M5").Value = -533.2941
$ws.Range("N5").Value = -1556.6666
$ws.Range("H32").Value = 2233.0894
$ws.Range("I32").Value = 2233.0894
$ws.Range("K32").Value = 2233.0894
$ws.Range("M32").Value = -1946.0894
$ws.Range("H45").Value = 1469.2
$ws.Range("I45").Value = 1198.8572
$ws.Range("J45").Value = 2100
$ws.Range("K45").Value = 1198.8572
$ws.Range("L45").Value = 2100
$ws.Range("M45").Value = -821.8571999999999
$ws.Range("N45").Value = -2854
$ws.Range("H61").Value = 250002620
$ws.Range("I61").Value = 333335500
$ws.Range("K61").Value = 333335500
$ws.Range("M61").Value = -333335288
$ws.Range("H74").Value = 52640020
$ws.Range("I74").Value = 76930650
$ws.Range("K74").Value = 76930650
$ws.Range("M74").Value = -76929776
$ws.Range("H77").Value = 52640020
$ws.Range("I77").Value = 76930650
$ws.Range("K77").Value = 384653250
$ws.Range("M77").Value = -384648882
$ws.Range("H97").Value = 523.6875
$ws.Range("I97").Value = 487.5
$ws.Range("K97").Value = 487.5
$ws.Range("M97").Value = 8.5
$ws.Range("H102").Value = 11112880
$ws.Range("I102").Value = 12501740
$ws.Range("K102").Value = 12501740
$ws.Range("M102").Value = -12500118
$ws.Range("H110").Value = 336670.66
$ws.Range("I110").Value = 502499.5
$ws.Range("J110").Value = 5013
$ws.Range("K110").Value = 502499.5
$ws.Range("L110").Value = 5013
$ws.Range("M110").Value = -500454.5
$ws.Range("N110").Value = -9103
$ws.Range("I116").Value = 921140.75
$ws.Range("J116").Value = 4124.75
$ws.Range("K116").Value = 921140.75
$ws.Range("L116").Value = 4124.75
$ws.Range("M116").Value = -918846.75
$ws.Range("N116").Value = -8712.75
$ws.Range("H132").Value = 4637667
$ws.Range("I132").Value = 2567573
$ws.Range("K132").Value = 7702719
$ws.Range("M132").Value = -7700189
$ws.Range("H136").Value = 250002620
$ws.Range("I136").Value = 333335500
$ws.Range("K136").Value = 1000006500
$ws.Range("M136").Value = -1000003950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 921140.75
$ws.Range("J3").Value = 4124.75
$ws.Range("K3").Value = 921140.75
$ws.Range("L3").Value = 4124.75
$ws.Range("M3").Value = -921026.75
$ws.Range("N3").Value = -4352.75
$ws.Range("H4").Value = 748.4
$ws.Range("I4").Value = 645.2941
$ws.Range("J4").Value = 1332.6666
$ws.Range("K4").Value = 645.2941
$ws.Range("L4").Value = 1332.6666
$ws.Range("M4").Value = -530.2941
$ws.Range("N4").Value = -1562.6666
$ws.Range("H20").Value = 1736.579
$ws.Range("I20").Value = 1579
$ws.Range("K20").Value = 1579
$ws.Range("M20").Value = -1332
$ws.Range("H86").Value = 3245.923
$ws.Range("I86").Value = 3144.5557
$ws.Range("K86").Value = 3144.5557
$ws.Range("M86").Value = -2021.5557
$ws.Range("H89").Value = 3245.923
$ws.Range("I89").Value = 3144.5557
$ws.Range("K89").Value = 15722.7785
$ws.Range("M89").Value = -10106.7785
$ws.Range("H105").Value = 3052.05
$ws.Range("I105").Value = 2483.1538
$ws.Range("K105").Value = 2483.1538
$ws.Range("M105").Value = -736.1538
$ws.Range("I134").Value = 19624346
$ws.Range("J134").Value = 1898.5
$ws.Range("K134").Value = 58873038
$ws.Range("L134").Value = 5695.5
$ws.Range("M134").Value = -58870503
$ws.Range("N134").Value = -10765.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3708.9207
$ws.Range("I31").Value = 2321.537
$ws.Range("K31").Value = 2321.537
$ws.Range("M31").Value = -2026.537
$ws.Range("H34").Value = 3708.9207
$ws.Range("I34").Value = 2321.537
$ws.Range("K34").Value = 2321.537
$ws.Range("M34").Value = -2119.537
$ws.Range("H50").Value = 5105
$ws.Range("I50").Value = 5105
$ws.Range("K50").Value = 5105
$ws.Range("M50").Value = -4480
$ws.Range("H60").Value = 18721.348
$ws.Range("I60").Value = 8773
$ws.Range("J60").Value = 20815.736
$ws.Range("K60").Value = 8773
$ws.Range("L60").Value = 20815.736
$ws.Range("M60").Value = -8262
$ws.Range("N60").Value = -21837.736
$ws.Range("H86").Value = 10964
$ws.Range("I86").Value = 8118.6665
$ws.Range("K86").Value = 8118.6665
$ws.Range("M86").Value = -6995.6665
$ws.Range("H89").Value = 10964
$ws.Range("I89").Value = 8118.6665
$ws.Range("K89").Value = 40593.3325
$ws.Range("M89").Value = -34977.3325
$ws.Range("H122").Value = 2971.48
$ws.Range("I122").Value = 2799.65
$ws.Range("K122").Value = 8398.950000000001
$ws.Range("M122").Value = -5948.950000000001
$ws.Range("H132").Value = 40001796
$ws.Range("I132").Value = 50001788
$ws.Range("K132").Value = 150005364
$ws.Range("M132").Value = -150002834
$ws.Range("H141").Value = 85359.60000000001
$ws.Range("J141").Value = 94288
$ws.Range("L141").Value = 94288
$ws.Range("N141").Value = -104648

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 881774.75
$ws.Range("I4").Value = 1158940.5
$ws.Range("K4").Value = 3476821.5
$ws.Range("M4").Value = -3476709.5
$ws.Range("H23").Value = 766
$ws.Range("I23").Value = 613.8
$ws.Range("K23").Value = 1841.4
$ws.Range("M23").Value = -1606.4
$ws.Range("H46").Value = 200
$ws.Range("J46").Value = 200
$ws.Range("L46").Value = 600
$ws.Range("N46").Value = -782
$ws.Range("H56").Value = 13672.808
$ws.Range("I56").Value = 13672.808
$ws.Range("K56").Value = 13672.808
$ws.Range("M56").Value = -13142.808
$ws.Range("H98").Value = 906.2222
$ws.Range("I98").Value = 906.2222
$ws.Range("K98").Value = 2718.6666
$ws.Range("M98").Value = -1220.6666
$ws.Range("H131").Value = 2130.6365
$ws.Range("I131").Value = 1929.75
$ws.Range("J131").Value = 2666.3333
$ws.Range("K131").Value = 5789.25
$ws.Range("L131").Value = 7998.999899999999
$ws.Range("M131").Value = -749.25
$ws.Range("N131").Value = -18078.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6898.5
$ws.Range("I70").Value = 6898.5
$ws.Range("K70").Value = 6898.5
$ws.Range("M70").Value = -6628.5
$ws.Range("H73").Value = 6898.5
$ws.Range("I73").Value = 6898.5
$ws.Range("K73").Value = 6898.5
$ws.Range("M73").Value = -5962.5
$ws.Range("H80").Value = 3039.5
$ws.Range("I80").Value = 3193
$ws.Range("J80").Value = 2579
$ws.Range("K80").Value = 3193
$ws.Range("L80").Value = 2579
$ws.Range("M80").Value = -2195
$ws.Range("N80").Value = -4575
$ws.Range("H83").Value = 3039.5
$ws.Range("I83").Value = 3193
$ws.Range("J83").Value = 2579
$ws.Range("K83").Value = 15965
$ws.Range("L83").Value = 12895
$ws.Range("M83").Value = -10973
$ws.Range("N83").Value = -22879
$ws.Range("H122").Value = 4872.6875
$ws.Range("I122").Value = 3039.762
$ws.Range("J122").Value = 8371.909
$ws.Range("K122").Value = 9119.286
$ws.Range("L122").Value = 25115.727
$ws.Range("M122").Value = -6669.286
$ws.Range("N122").Value = -30015.727

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2839.6428
$ws.Range("I22").Value = 2797.111
$ws.Range("K22").Value = 2797.111
$ws.Range("M22").Value = -2502.111
$ws.Range("H27").Value = 2839.6428
$ws.Range("I27").Value = 2797.111
$ws.Range("K27").Value = 2797.111
$ws.Range("M27").Value = -2690.111
$ws.Range("H46").Value = 1647.3
$ws.Range("I46").Value = 1728
$ws.Range("J46").Value = 1459
$ws.Range("K46").Value = 1728
$ws.Range("L46").Value = 1459
$ws.Range("M46").Value = -1540
$ws.Range("N46").Value = -1835
$ws.Range("H47").Value = 13333.333
$ws.Range("H52").Value = 13333.333
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H61").Value = 7009.5557
$ws.Range("I61").Value = 7260.75
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 7260.75
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -7058.75
$ws.Range("N61").Value = -5404
$ws.Range("H82").Value = 2810.6667
$ws.Range("I82").Value = 2500
$ws.Range("J82").Value = 2966
$ws.Range("K82").Value = 2500
$ws.Range("L82").Value = 2966
$ws.Range("M82").Value = -2139
$ws.Range("N82").Value = -3688
$ws.Range("H85").Value = 2810.6667
$ws.Range("I85").Value = 2500
$ws.Range("J85").Value = 2966
$ws.Range("K85").Value = 2500
$ws.Range("L85").Value = 2966
$ws.Range("M85").Value = -1252
$ws.Range("N85").Value = -5462
$ws.Range("H93").Value = 1064.3334
$ws.Range("I93").Value = 1071.3158
$ws.Range("J93").Value = 998
$ws.Range("K93").Value = 1071.3158
$ws.Range("L93").Value = 998
$ws.Range("M93").Value = 176.6841999999999
$ws.Range("N93").Value = -3494
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H113").Value = 7009.5557
$ws.Range("I113").Value = 7260.75
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 7260.75
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -5090.75
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 2408
$ws.Range("I136").Value = 2489.6428
$ws.Range("K136").Value = 7468.928400000001
$ws.Range("M136").Value = -4918.928400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4889.4375
$ws.Range("I62").Value = 2985.4
$ws.Range("K62").Value = 2985.4
$ws.Range("M62").Value = -2361.4
$ws.Range("H65").Value = 4889.4375
$ws.Range("I65").Value = 2985.4
$ws.Range("K65").Value = 14927
$ws.Range("M65").Value = -11807
$ws.Range("H122").Value = 1643.6
$ws.Range("I122").Value = 1579.3
$ws.Range("J122").Value = 1772.2
$ws.Range("K122").Value = 4737.9
$ws.Range("L122").Value = 5316.6
$ws.Range("M122").Value = -2287.9
$ws.Range("N122").Value = -10216.6
$ws.Range("H136").Value = 29414188
$ws.Range("I136").Value = 31252356
$ws.Range("K136").Value = 93757068
$ws.Range("M136").Value = -93754518
$ws.Range("H141").Value = 83196.8
$ws.Range("J141").Value = 79496
$ws.Range("L141").Value = 79496
$ws.Range("N141").Value = -89856
